$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: set B/C/E via .Value directly (never numeric-looking),
# and D via NumberFormat="@" + .Value + .Style reset so it stays text like the source.

$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.232.70'
$ws.Range("D2").Style = "Normal"

$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.882.13'
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = '  +0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("D4").Style = "Normal"

$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.25'
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5139'
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = '  +1.40%  '

$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08386'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  +0.67%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.73'
$ws.Range("D11").Style = "Normal"

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.260'
$ws.Range("D12").Style = "Normal"

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.63'
$ws.Range("D13").Style = "Normal"

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.866.06'
$ws.Range("D14").Style = "Normal"

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.276'
$ws.Range("D15").Style = "Normal"

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.010'
$ws.Range("D16").Style = "Normal"

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("D17").Style = "Normal"

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.07'
$ws.Range("D18").Style = "Normal"

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06715'
$ws.Range("D19").Style = "Normal"

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.80'
$ws.Range("D20").Style = "Normal"

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.008'
$ws.Range("D21").Style = "Normal"

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.031'
$ws.Range("D22").Style = "Normal"

$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.276.36'
$ws.Range("D23").Style = "Normal"

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("D24").Style = "Normal"

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.266'
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = '  +1.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.70'
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = '  -2.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.462'
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.69'
$ws.Range("D28").Style = "Normal"

$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '125.62'
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1057'
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.874'
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.624'
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = '  +1.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.576'
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02449'
$ws.Range("D35").Style = "Normal"

$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06576'
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = '  +1.92%  '

$ws.Range("E38").Value = '  -0.49%  '

$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6487'
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = '  +1.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.244'
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.004'
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = '  +1.03%  '

$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6099'
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.08'
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = '  +1.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.699'
$ws.Range("D45").Style = "Normal"

$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.280'
$ws.Range("D46").Style = "Normal"

$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.014'
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = '  +1.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.234'
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.14'
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06925'
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = '  -0.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.89'
$ws.Range("D51").Style = "Normal"
